# The sheet originally has columns:
#   A RegistryCode | B age | C WSO | D FPM | E WPM | F WPM_log | G FPM_log
#
# The commit "create .xlsx-df with only log rt" drops the raw-rate columns
# (FPM, WPM) and keeps only the log-transformed rates, which slide left
# into D/E:
#   A RegistryCode | B age | C WSO | D WPM_log | E FPM_log
#
# Deleting the entire column D twice removes FPM first (shifting WPM,
# WPM_log, FPM_log left by one) and then removes the now-shifted WPM
# (shifting WPM_log, FPM_log left by one more), leaving WPM_log in D and
# FPM_log in E exactly as in the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(4).Delete() | Out-Null
$ws.Columns.Item(4).Delete() | Out-Null
